# "added pascal and jacobi-1D vary blocksize experiments"
#
# Target layout (rows 70-78, columns A:B) after the edit:
#   70: A=1  B="Input size"
#   71: A=2  B="Quantitative Comparison to Related Work - looking for other methods that apply to cyclic or block cyclic "
#   72: A=3  B="Writing Improvements"
#   73: A=4  B="Slowdowns for existing benchmarks"
#   74: A=5  B="Varying the block size"            <- newly inserted row
#   75: (blank)
#   76: B="Finish 1 to 3 by September 26"
#   77: B="Finish 4 by Oct 15"
#   78: B="If we can do case study thing, do it by Sept 26"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row 74 (pushes the blank separator row + the three
# "Finish ..." rows down by one, from 75/76/77 to 76/77/78).
$ws.Rows(74).Insert()

# Populate the newly inserted row with the "Varying the block size" item.
$ws.Range("A74").Value = 5
$ws.Range("B74").Value = "Varying the block size"

# Rename the "Comparison to Related Work" row (B71) into the new, more
# specific quantitative-comparison description.
$ws.Range("B71").Value = "Quantitative Comparison to Related Work - looking for other methods that apply to cyclic or block cyclic "

# Move the selection to match the new bottom of the list.
$ws.Range("B80").Select() | Out-Null
